# Add two new ship-type rows (110-er and 135-er) to the fleet database.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76: 110-er
$ws.Range("A76").Value = "110-er"
$ws.Range("B76").Value = 11.5
$ws.Range("C76").Value = 110
$ws.Range("D76").Value = 3000
$ws.Range("E76").Value = "M8"

# Row 77: 135-er
$ws.Range("A77").Value = "135-er"
$ws.Range("B77").Value = 11.5
$ws.Range("C77").Value = 135
$ws.Range("D77").Value = 4000
$ws.Range("E77").Value = "M9"

# Match column A style (left-aligned) used by the other data rows.
$ws.Range("A76:A77").HorizontalAlignment = -4131

# Select the first empty row beneath the table, matching the author's
# final view/selection.
$ws.Range("A78").Select() | Out-Null
